# Apply the "Added Prospectus supp figs" edit:
#  - clear the "Significant" (X) marker in column J for rows 4, 7 and 24
#  - highlight rows 4 and 7 with a new custom pink fill color
#  - recolor row 24 to match the existing "green" highlight (same as row 3, 12, ...)
#  - restore the view/selection to match the saved state (top of sheet, F11 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the Significant flags that are no longer set ---
$ws.Range("J4").Clear()
$ws.Range("J7").Clear()
$ws.Range("J24").Clear()

# --- recolor rows 4 and 7 with a new custom fill (pink, RGB FF66CC) ---
$ws.Range("A4:I4").Interior.Color = 13395711
$ws.Range("A7:I7").Interior.Color = 13395711

# --- recolor row 24 to the same green fill already used elsewhere (row 3, 12, ...) ---
$ws.Range("A3:I3").Copy()
$ws.Range("A24:I24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- restore view state: scroll back to top, select F11 ---
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("F11").Select()

Write-Output "edit applied"
